# Generate Report for Handback
# This script updates the localization-status workbook's zh-cn and de-de
# sheets with the results of the handback report run: it records the
# Latest Target File / Latest Handback File / Latest Handback DateTime for
# the 618f9b18 row, and records the error detail for a stale handback.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b87dfe1efeb76fc789c02be17ce20cc4c59721c2/e2e/618f9b18-a87f-4a8e-be60-112c85bb74ef.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b4a325e397685a5b62c2e2e9c9ee05e76902701/e2e/618f9b18-a87f-4a8e-be60-112c85bb74ef.md."

# Magic input for ColumnWidth that round-trips to a raw OOXML column
# <col width=.../> of exactly 40 in this engine.
$fortyColWidth = 39.166666666666664

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P, the 16th column) to fit the message.
$wsZh.Columns.Item(16).ColumnWidth = $fortyColWidth

# Latest Target File (I5): now known, becomes a hyperlink.
$wsZh.Range("I5").Value = "618f9b18-a87f-4a8e-be60-112c85bb74ef.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fd4f0e01296ee5c592c0453312e71e0a93768378/e2e/618f9b18-a87f-4a8e-be60-112c85bb74ef.md", "", "", "618f9b18-a87f-4a8e-be60-112c85bb74ef.md") | Out-Null
$wsZh.Range("I5").Style = "HyperLink"
$wsZh.Range("I5").Font.Underline = 2
$wsZh.Range("I5").Font.Color = 15570276

# Latest Handback File (J5).
$wsZh.Range("J5").Value = "618f9b18-a87f-4a8e-be60-112c85bb74ef.fa417c05604ef6edca4a782bbde7a2195d31de9b.zh-cn.xlf"

# Latest Handback DateTime (K5).
$wsZh.Range("K5").Value = "2016-10-21 04:02:59"

# Error Detail (P5).
$wsZh.Range("P5").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P, the 16th column) to fit the message.
$wsDe.Columns.Item(16).ColumnWidth = $fortyColWidth

# Latest Target File (I5): now known, becomes a hyperlink.
$wsDe.Range("I5").Value = "618f9b18-a87f-4a8e-be60-112c85bb74ef.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/28a0a9a61dc2f9c1229a6e623a4cc95eb420a6d8/e2e/618f9b18-a87f-4a8e-be60-112c85bb74ef.md", "", "", "618f9b18-a87f-4a8e-be60-112c85bb74ef.md") | Out-Null
$wsDe.Range("I5").Style = "HyperLink"
$wsDe.Range("I5").Font.Underline = 2
$wsDe.Range("I5").Font.Color = 15570276

# Latest Handback File (J5).
$wsDe.Range("J5").Value = "618f9b18-a87f-4a8e-be60-112c85bb74ef.fa417c05604ef6edca4a782bbde7a2195d31de9b.de-de.xlf"

# Latest Handback DateTime (K5).
$wsDe.Range("K5").Value = "2016-10-21 04:03:18"

# Error Detail (P5).
$wsDe.Range("P5").Value = $errorDetail
